# Update the handback/handoff timestamps to reflect the newly generated
# report timestamps (commit: "Generate Report for Handback").

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-05 13:21:16"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-05 13:21:08"
$wsZhCn.Range("K2").Value = "2016-09-05 13:21:55"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-05 13:22:09"
